$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values per repull/recalculation of data
$ws.Range("F3").Value = -14
$ws.Range("F4").Value = -4
$ws.Range("F5").Value = -16
$ws.Range("F6").Value = -1
$ws.Range("F10").Value = -7
$ws.Range("F11").Value = -1
